{"js": "// Insert a new \"BloodType:\" paragraph right after the Stock Table's\n// \"...CenterID:<br/>BloodGroupID:<br/>Quantity:\" paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the target paragraph robustly: the Stock Table entry paragraph\n// that contains \"BloodGroupID:\" and ends with \"Quantity:\" (there is another,\n// unrelated \"Quantity:\" run earlier in the Advance Request Table, but it is\n// not the last run of its paragraph, so this combined check is unambiguous).\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"BloodGroupID:\") !== -1 && /Quantity:$/.test(text)) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Could not find the \"Quantity:\" paragraph of the Stock Table.');\n}\n\n// Insert the new paragraph right after it; it naturally inherits the\n// sz/szCs (24) run formatting from the paragraph it follows.\ntarget.insertParagraph(\"BloodType:\", \"After\");\n\nawait context.sync();\n", "ps1": "# Insert a new \"BloodType:\" paragraph right after the Stock Table's\n# \"...CenterID:<br/>BloodGroupID:<br/>Quantity:\" paragraph.\n$d = $word.ActiveDocument\n\n# Locate the target paragraph robustly: the Stock Table entry paragraph\n# that contains \"BloodGroupID:\" and ends with \"Quantity:\" (there is another,\n# unrelated \"Quantity:\" run earlier in the Advance Request Table, but it is\n# not the last run of its paragraph, so this combined check is unambiguous).\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*BloodGroupID:*\" -and $t -match \"Quantity:\\s*$\") {\n        $targetIndex = $i\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not find the 'Quantity:' paragraph of the Stock Table.\"\n}\n\n$target = $d.Paragraphs.Item($targetIndex)\n$null = $target.Range.InsertParagraphAfter()\n\n# The handle returned by InsertParagraphAfter is stale for further edits, so\n# re-fetch the freshly-inserted paragraph from the collection by index and\n# set its text. It naturally inherits the sz/szCs (24) run formatting from\n# the paragraph it follows.\n$newPara = $d.Paragraphs.Item($targetIndex + 1)\n$newPara.Range.Text = \"BloodType:\"\n"}
